$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before the old "Q" column (most_frequent_value),
# shifting the existing Q,R,S,T columns to S,T,U,V.
$ws.Range("Q1:Q8").EntireColumn.Insert()
$ws.Range("Q1:Q8").EntireColumn.Insert()

# New header cells for the two inserted columns.
$ws.Range("Q1").Value = "default_count"
$ws.Range("R1").Value = "default_value"

# New per-row data for the inserted columns.
$ws.Range("Q2").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("Q8").Value = 0

$ws.Range("R2").Value = "<Unspecified>"
$ws.Range("R3").Value = "<Unspecified>"
$ws.Range("R4").Value = "<Unspecified>"
$ws.Range("R5").Value = "<Unspecified>"
$ws.Range("R6").Value = "<Unspecified>"
$ws.Range("R7").Value = "<Unspecified>"
$ws.Range("R8").Value = "<Unspecified>"

# Update the (now shifted) most_frequent_value column S with new values.
# Cells whose new text looks like a number (S2, S4, S6) need to be forced to
# Text format first so Excel keeps them as strings instead of auto-converting.
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "1"
$ws.Range("S3").Value = "University"
$ws.Range("S4").NumberFormat = "@"
$ws.Range("S4").Value = "24029000"
$ws.Range("S5").Value = 'Collge d''enseignement gZnZral et professionnel (CfGEP) de Thetford'
$ws.Range("S6").NumberFormat = "@"
$ws.Range("S6").Value = "51.2205"
$ws.Range("S7").Value = "Ecology"
$ws.Range("S8").Value = "no"
